$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("10__bc529b88-4ffa-36")
$ws.Range("A5").Value = "10 Inventories nan"
$ws.Range("A9").Value = "10 Inventories nan"

$ws = $wb.Worksheets.Item("12__7f5b7384-15a4-38")
$ws.Range("A2").Value = "None Trade payables"
$ws.Range("A3").Value = "None Amounts payable to controlling entity (Refer Note 19)"
$ws.Range("A4").Value = "None Other creditors and accruals"
$ws.Range("A5").Value = "None nan"
$ws.Range("A6").Value = "None Trade payables"
$ws.Range("A7").Value = "None Amounts payable to controlling entity (Refer Note 19)"
$ws.Range("A8").Value = "None Other creditors and accruals"
$ws.Range("A9").Value = "None nan"

$ws = $wb.Worksheets.Item("14__7f5b7384-15a4-38")
$ws.Range("A3").Value = "Current nan"
$ws.Range("A5").Value = "Non-current nan"
$ws.Range("A7").Value = "Current nan"
$ws.Range("A9").Value = "Non-current nan"

$ws = $wb.Worksheets.Item("15__c932dc91-8c71-3d")
$ws.Range("A2").Value = "None at 1 January"
$ws.Range("A3").Value = "None at31 December"
$ws.Range("A4").Value = "None nan"
$ws.Range("A5").Value = "None nan"
$ws.Range("A6").Value = "None nan"
$ws.Range("A7").Value = "None at 1 January"
$ws.Range("A8").Value = "None at 31 December"
$ws.Range("A9").Value = "None at 1 January"
$ws.Range("A10").Value = "None at31 December"
$ws.Range("A11").Value = "None nan"
$ws.Range("A12").Value = "None nan"
$ws.Range("A13").Value = "None nan"
$ws.Range("A14").Value = "None at 1 January"
$ws.Range("A15").Value = "None at 31 December"
$ws.Range("A16").Value = "None at 1 January"
$ws.Range("A17").Value = "None at31 December"
$ws.Range("A18").Value = "None nan"
$ws.Range("A19").Value = "None nan"
$ws.Range("A20").Value = "None nan"
$ws.Range("A21").Value = "None at 1 January"
$ws.Range("A22").Value = "None at 31 December"
$ws.Range("A23").Value = "None at 1 January"
$ws.Range("A24").Value = "None at31 December"
$ws.Range("A25").Value = "None nan"
$ws.Range("A26").Value = "None nan"
$ws.Range("A27").Value = "None nan"
$ws.Range("A28").Value = "None at 1 January"
$ws.Range("A29").Value = "None at 31 December"

$ws = $wb.Worksheets.Item("15__5213ad68-9859-37")
$ws.Range("A2").Value = "None On issue at 1 January"
$ws.Range("A3").Value = "None On issue at 31 December"
$ws.Range("A4").Value = "None On issue at 1 January"
$ws.Range("A5").Value = "None On issue at 31 December"

$ws = $wb.Worksheets.Item("15__5e4f56b4-95fd-32")
$ws.Range("A2").Value = "None Cents per share (fully franked)"
$ws.Range("A3").Value = "None Total dollar amount"
$ws.Range("A4").Value = "None Cents per share (fully franked)"
$ws.Range("A5").Value = "None Total dollar amount"

